$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 held product 3841 "Fringe Poncho". Replace it in place with the
# new product 3842 "Luxe Turtleneck" (same CareInstr text in column G).
$a2 = $ws.Cells.Item(2, 1)
$a2.NumberFormat = "@"
$a2.Value = "3842"
$a2.ClearFormats()

$ws.Cells.Item(2, 2).Value = "Luxe Turtleneck"
$ws.Cells.Item(2, 3).Value = "Winter White"
$ws.Cells.Item(2, 4).Value = 99
$ws.Cells.Item(2, 5).Value = "XS-XL"
$ws.Cells.Item(2, 6).Value = "51% Cotton, 30% Polyester, 19% Nylon"
